# Fix the "Samples" tab query (cell B3 on the "startup" sheet): the Tumor
# column previously used the aggregated/collected `tumor` list variable,
# which is wrong - it should read the sample's own sample_tumor_status field.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newTumorQuery = "MATCH (s:study)<--(p:participant)<--(samp:sample)`r`nWHERE s.study_name in [`"Childhood Cancer Data Initiative (CCDI): Free the Data: Open Sharing of Comprehensive Genomic Childhood Cancer Datasets (Kansas)`"]`r`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`r`nRETURN  `r`n coalesce(samp.sample_id, '') as ``Sample ID``,`r`n coalesce(p.participant_id,'') as ``Participant ID``,`r`n coalesce(s.study_name, '') as ``Study Name``,`r`n coalesce(s.phs_accession,'') as ``Accession``,`r`ncoalesce(samp.sample_tumor_status,'') as ``Tumor``,`r`ncoalesce(samp.sample_type,'') as ``Analyte Type```r`nORDER By samp.sample_id LIMIT 100"

$ws.Range("B3").Value = $newTumorQuery

# The author's last selection in the sheet ended up on B4 (FilesTab query cell)
$ws.Range("B4").Select()
